# Rename the "Collection_MH" tab to "CRF_MH" (per commit: "rename Collection to CRF in tabs").
# Renaming via the Worksheet.Name property also keeps any dependent references
# (e.g. the sheet-scoped _xlnm._FilterDatabase defined name) pointing at the
# correct sheet, since Excel updates those automatically on rename.
$wb = $excel.ActiveWorkbook

$oldName = "Collection_MH"
$newName = "CRF_MH"

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq $oldName) {
        $ws = $sheet
        break
    }
}

if ($ws -eq $null) {
    # Fallback: operate on the active sheet if the expected name wasn't found.
    $ws = $wb.ActiveSheet
}

$ws.Name = $newName
